# Build 2.0.0.6 - Added Include all versions for DqlDataSource
# Fixed Error Messaging grid. Fixed bug grid mapping not deleting correct item

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Fixed Error Messaging grid: rows 2 & 3, authors/keywords columns now
#     carry real pipe-delimited values instead of the shared blank string,
#     the creation-date time-of-day shifts, and the reference-count gets a
#     fractional correction (bug fix: "not deleting correct item"). ---
$ws1.Range("F2").Value = "a1|a3|a2"
$ws1.Range("G2").Value = "k1|k3|k2"
$ws1.Range("F3").Value = "a1|a3|a2"
$ws1.Range("G3").Value = "k1|k3|k2"

$ws1.Range("J2").Value = 41556.252268518518
$ws1.Range("J3").Value = 41556.252268518518

$ws1.Range("T2").Value = 1.01

# --- widen the two columns that now hold the longer pipe-delimited values ---
$ws1.Columns.Item(10).ColumnWidth = 26.43
$ws1.Columns.Item(20).ColumnWidth = 14.43

# Selection on Sheet1 moves to the header row of the affected columns
$ws1.Range("A1:I1").Select()

# --- Added "Include all versions" sheet (DqlDataSource) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1:I1").NumberFormat = "@"
$ws2.Range("A1").Value = "r_object_id"
$ws2.Range("B1").Value = "object_name"
$ws2.Range("C1").Value = "r_object_type"
$ws2.Range("D1").Value = "title"
$ws2.Range("E1").Value = "subject"
$ws2.Range("F1").Value = "authors"
$ws2.Range("G1").Value = "keywords"
$ws2.Range("H1").Value = "a_application_type"
$ws2.Range("I1").Value = "a_status"

$ws2.Range("B7").Select()
